$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("promotions")

$ws.Range("A2").Value = "REACH Tuition Course Fees 2021.pdf"
$ws.Range("B2").Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/onshore/coe/reach/Reach_Tuition_Course_Fees_2021_v1.0.pdf"

$ws.Range("A3").Value = "REACH Q4 Promotions.pdf"
$ws.Range("B3").Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/onshore/coe/reach/Reach-NonCoE_CoE-Q4-Promotions-1OCT-31DEC21_Vol-1.0.pdf"

$ws.Range("B3").Select()

$wb.Save()
